$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Four new product rows appended below the existing data (rows 2-5).
# Column C ("Bar Code") holds numeric-looking values that must be stored
# as text (shared strings), matching the existing rows' cell type -
# prefixing with an apostrophe forces Excel to treat the entry as text.

$ws.Range("A6").Value = "Cisflem (Carbo) 125 mg/60 ml Syrup"
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = "'7"
$ws.Range("D6").Value = $true

$ws.Range("A7").Value = "(Amoxicillin) 125 mg/60 ml Syrup"
$ws.Range("B7").Value = 50
$ws.Range("C7").Value = "'8"
$ws.Range("D7").Value = $true

$ws.Range("A8").Value = "Mucosolve (Ambroxol) 30 mg/60 ml Syrup"
$ws.Range("B8").Value = 50
$ws.Range("C8").Value = "'9"
$ws.Range("D8").Value = $true

$ws.Range("A9").Value = "Cisflem (Carbo) 125 mg/60 ml Syrup"
$ws.Range("B9").Value = 50
$ws.Range("C9").Value = "'10"
$ws.Range("D9").Value = $true

# Match the formatting of the rest of the data rows (wrap-text, default font).
$ws.Range("A6:D9").Style = "Normal"
$ws.Range("A6:D9").WrapText = $true
